$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 6

# Copy formatting (incl. the date number format style) from the cell above
# so the new row's date cell reuses the existing style instead of creating
# a brand-new one.
$ws.Cells.Item(2, 1).Copy($ws.Cells.Item($row, 1))
$ws.Cells.Item($row, 1).Value = 42588.471516203703

$ws.Cells.Item($row, 2).Value = "Bag"
$ws.Cells.Item($row, 3).Value = 6541
$ws.Cells.Item($row, 4).Value = 10281
$ws.Cells.Item($row, 5).Value = 1264
$ws.Cells.Item($row, 6).Value = 121
$ws.Cells.Item($row, 7).Value = 69
$ws.Cells.Item($row, 8).Value = 62
$ws.Cells.Item($row, 9).Value = 35
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 2
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = 100
